$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.333.55"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.626.78"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.12"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3740"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3625"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.25"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08156"
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.222"
$ws.Range("E11").Value = "  -2.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.22"
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.482"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001242"
$ws.Range("E15").Value = "  -2.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.298"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "1.625.92"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.83"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06958"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.50"
$ws.Range("E20").Value = "  -3.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.509"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.56"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").Value = "23.335.86"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.110"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.463"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.28"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.60"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.296"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.19"
$ws.Range("E30").Value = "  -2.59%  "
$ws.Range("D31").Value = "1.806.06"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.235"
$ws.Range("E32").Value = "  -3.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.742"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.031"
$ws.Range("E34").Value = "  +7.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.68"
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02754"
$ws.Range("E36").Value = "  -2.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2503"
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08768"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07090"
$ws.Range("E39").Value = "  -2.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.941"
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.340"
$ws.Range("E41").Value = "  -2.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6980"
$ws.Range("E42").Value = "  -1.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.00"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.07"
$ws.Range("E44").Value = "  -3.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6500"
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.272"
$ws.Range("E47").Value = "  -2.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.968"
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07979"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.191"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.52"
$ws.Range("E51").Value = "  -2.42%  "
